$d = $word.ActiveDocument

# NOTE 1: Paragraph/Range handles captured *before* a structural edit (one
# that adds/removes paragraphs) are not safe to keep using afterwards in
# this host - they are resolved against the live document by numeric
# index, so a stale handle can silently point at the wrong paragraph once
# the paragraph count has shifted. Every step below therefore re-resolves
# positions fresh (via Find, anchored off stable nearby text) immediately
# before each mutation.
#
# NOTE 2: the document already has one "_GoBack" bookmark. Bookmark names
# must stay unique, so the existing one is removed *before* the new one
# (at its relocated position) is added - doing it the other way round
# leaves two same-named bookmarks alive at once and a later
# Bookmarks("_GoBack") lookup/delete becomes ambiguous.

function Remove-CommentsParagraphAfter($anchorText) {
    # Finds the paragraph containing $anchorText, then deletes the very
    # next paragraph (the "Comments:" label paragraph that immediately
    # follows it) by merging it away into the paragraph after that.
    # Returns the (post-delete) index of the paragraph that absorbed it.
    $rng = $d.Content
    $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $anchorIdx = $rng.Paragraphs(1).Index
    $commentsIdx = $anchorIdx + 1

    $startPos = $d.Paragraphs($commentsIdx).Range.Start
    $endPos = $d.Paragraphs($commentsIdx + 1).Range.Start
    $d.Range($startPos, $endPos).Delete()

    return $commentsIdx
}

# --- A) Drop the existing "_GoBack" bookmark near the end of the document
#        and fold the two runs around it back into a single run. ----------
$d.Bookmarks("_GoBack").Delete()
$d.Content.Find.Execute(" (for the weighted average model).", $true, $false, $false, $false, $false, $true, 1, $false, " (for the weighted average model).", 2) | Out-Null

# --- B) First "Comments:" block (Stage 1 section): the "Comments:"
#        paragraph is merged away into the following (empty) paragraph,
#        which ends up holding the relocated "_GoBack" bookmark. ---------
$mergedIdx1 = Remove-CommentsParagraphAfter "For each ortholog pair- how much data do we have for it?"
$d.Bookmarks.Add("_GoBack", $d.Paragraphs($mergedIdx1).Range) | Out-Null

# --- C) Fourth "Comments:" block (Stage 4 section): the "Comments:"
#        paragraph is removed outright, leaving the following (already
#        empty) bold/underlined paragraph untouched. ----------------------
Remove-CommentsParagraphAfter "for each ortholog couple perform lasso regression 100 times (bootstrapping)." | Out-Null
